$d = $word.ActiveDocument
Write-Output $d.list_commands()
